# s2cDNASample (H.BROWN, 10.29.19) - "further cleaning to metadata"
#
#  1. The cDNA sample number (column G, shared by every data row) is
#     corrected from E7760 to E7420.
#  2. Column G (s2cDNASampleNumber) is restyled onto its own dedicated
#     cell style - General number format, Arial 11 - instead of reusing
#     whatever style each row happened to inherit.
#  3. Column H (roboticS2Prep) is rewritten as an explicit =FALSE()
#     formula on every row instead of a bare boolean literal.
#  4. The sheet's remembered selection moves from the H column to the G
#     column (G2:G41 / active cell G2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 41
$colG = 7   # s2cDNASampleNumber
$colH = 8   # roboticS2Prep

# 1. Fix the sample number value for every data row (it is the same
#    shared string for all of them, so this rewrites that single shared
#    string in place).
$ws.Range("G$firstRow`:G$lastRow").Value = "E7420"

# 2 & 3. Re-style column G and turn column H into a live =FALSE() formula,
#    one row at a time so each cell gets its own (non-shared) formula,
#    matching how the sheet originally stored it.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $gCell = $ws.Cells.Item($r, $colG)
    $gCell.Font.Size = 11
    $gCell.Font.Name = "Arial"

    $hCell = $ws.Cells.Item($r, $colH)
    $hCell.Formula = "=FALSE()"
}

# 4. Move the sheet's active selection from H2:H41 to G2:G41.
$ws.Range("G$firstRow`:G$lastRow").Select()
